# Generate Report for Handoff
# Inserts a new tracked file ("b013591f-7584-4a73-a059-a39a2aac7c99.md") into the
# localization-status workbook, one row above the existing last row (which held
# "b04048be-ba5b-487b-a279-7de58c31591f.md") on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Push the current last data row (row 7, "b04048be...") down to row 8 and open
# up a fresh row 7 for the new entry.
$ws1.Rows.Item(7).Insert()

$ws1.Cells.Item(7, 1).Value = "b013591f-7584-4a73-a059-a39a2aac7c99.md"
$ws1.Cells.Item(7, 2).Value = "e2e\b013591f-7584-4a73-a059-a39a2aac7c99.md"
$ws1.Cells.Item(7, 3).Value = ".md"
$ws1.Cells.Item(7, 4).Value = ""
$ws1.Cells.Item(7, 5).Value = "Ready for handoff"
$ws1.Cells.Item(7, 6).Value = "Ready for handoff"
$ws1.Cells.Item(7, 7).Value = "2017-02-09 08:01:57"
$ws1.Cells.Item(7, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild every hyperlink on the sheet (row-insert does not itself relocate the
# hyperlink objects), now spanning rows 2-8.
$ws1.Cells.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Cells.Item(2, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a28272383e24ad73c8d92fe07e2d8956f1af8d3/e2e/17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md", "", "", "e2e\17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(3, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692db1b7948673fd1feb1023296c1e0e1ec9444a/e2e/23ea27e7-c26e-4c26-aa4c-af80ea04e4dc.md", "", "", "e2e\23ea27e7-c26e-4c26-aa4c-af80ea04e4dc.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(4, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1af3634fd0a2a83d98831a73295a81a69f4d53d3/e2e/38b81aa5-ef24-46c7-a1ab-2ffc89f9234c.md", "", "", "e2e\38b81aa5-ef24-46c7-a1ab-2ffc89f9234c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(5, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/108c5fd2cde90b86983be5bb4dee58dfd974048e/e2e/85e3d7a9-c5e2-4c10-8003-d5d93f2db402.md", "", "", "e2e\85e3d7a9-c5e2-4c10-8003-d5d93f2db402.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(6, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692db1b7948673fd1feb1023296c1e0e1ec9444a/e2e/a7341649-f68e-4d0f-911b-33b6545d0ef9.md", "", "", "e2e\a7341649-f68e-4d0f-911b-33b6545d0ef9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(7, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/408e5bacc2e22a8dd0b7d42bca0f848975663a38/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", "e2e\b013591f-7584-4a73-a059-a39a2aac7c99.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(8, 2), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/decaa4d912507b61ae7fd3fd8f9bb69ddc6d7150/e2e/b04048be-ba5b-487b-a279-7de58c31591f.md", "", "", "e2e\b04048be-ba5b-487b-a279-7de58c31591f.md") | Out-Null

# Grow the "Overview" table so the new row is part of it.
$lo1 = $ws1.ListObjects.Item("Overview")
$lo1.Resize($ws1.Range("A1:G8"))

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(7).Insert()

$ws2.Cells.Item(7, 1).Value = "b013591f-7584-4a73-a059-a39a2aac7c99.md"
$ws2.Cells.Item(7, 2).Value = ".md"
$ws2.Cells.Item(7, 3).Value = "Ready for handoff"
$ws2.Cells.Item(7, 4).Value = "e2e"
$ws2.Cells.Item(7, 5).Value = "ht"
$ws2.Cells.Item(7, 6).Value = "'False"
$ws2.Cells.Item(7, 7).Value = "b013591f-7584-4a73-a059-a39a2aac7c99.408e5bacc2e22a8dd0b7d42bca0f848975663a38.zh-cn.xlf"
$ws2.Cells.Item(7, 8).Value = "2017-02-09 08:01:40"
$ws2.Cells.Item(7, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(7, 9).Value = ""
$ws2.Cells.Item(7, 10).Value = ""
$ws2.Cells.Item(7, 11).Value = ""
$ws2.Cells.Item(7, 12).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(7, 12).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Cells.Item(7, 13).Value = ""
$ws2.Cells.Item(7, 14).Value = ""
$ws2.Cells.Item(7, 15).Value = "'True"
$ws2.Cells.Item(7, 16).Value = ""
$ws2.Cells.Item(7, 17).Value = "'False"
$ws2.Cells.Item(7, 18).Value = ""

# Rebuild hyperlinks: A2..A8 (source file links) plus the untouched J2 link.
$ws2.Cells.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a28272383e24ad73c8d92fe07e2d8956f1af8d3/e2e/17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md", "", "", "17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2, 10), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/4548cc233b7534afa9f4d73e7b6839a56d1ddc01/e2e/17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md", "", "", "17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692db1b7948673fd1feb1023296c1e0e1ec9444a/e2e/23ea27e7-c26e-4c26-aa4c-af80ea04e4dc.md", "", "", "23ea27e7-c26e-4c26-aa4c-af80ea04e4dc.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1af3634fd0a2a83d98831a73295a81a69f4d53d3/e2e/38b81aa5-ef24-46c7-a1ab-2ffc89f9234c.md", "", "", "38b81aa5-ef24-46c7-a1ab-2ffc89f9234c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/108c5fd2cde90b86983be5bb4dee58dfd974048e/e2e/85e3d7a9-c5e2-4c10-8003-d5d93f2db402.md", "", "", "85e3d7a9-c5e2-4c10-8003-d5d93f2db402.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(6, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692db1b7948673fd1feb1023296c1e0e1ec9444a/e2e/a7341649-f68e-4d0f-911b-33b6545d0ef9.md", "", "", "a7341649-f68e-4d0f-911b-33b6545d0ef9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(7, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/408e5bacc2e22a8dd0b7d42bca0f848975663a38/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", "b013591f-7584-4a73-a059-a39a2aac7c99.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(8, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/decaa4d912507b61ae7fd3fd8f9bb69ddc6d7150/e2e/b04048be-ba5b-487b-a279-7de58c31591f.md", "", "", "b04048be-ba5b-487b-a279-7de58c31591f.md") | Out-Null

$lo2 = $ws2.ListObjects.Item("zh-cn")
$lo2.Resize($ws2.Range("A1:R8"))

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(7).Insert()

$ws3.Cells.Item(7, 1).Value = "b013591f-7584-4a73-a059-a39a2aac7c99.md"
$ws3.Cells.Item(7, 2).Value = ".md"
$ws3.Cells.Item(7, 3).Value = "Ready for handoff"
$ws3.Cells.Item(7, 4).Value = "e2e"
$ws3.Cells.Item(7, 5).Value = "ht"
$ws3.Cells.Item(7, 6).Value = "'False"
$ws3.Cells.Item(7, 7).Value = "b013591f-7584-4a73-a059-a39a2aac7c99.408e5bacc2e22a8dd0b7d42bca0f848975663a38.de-de.xlf"
$ws3.Cells.Item(7, 8).Value = "2017-02-09 08:01:57"
$ws3.Cells.Item(7, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(7, 9).Value = ""
$ws3.Cells.Item(7, 10).Value = ""
$ws3.Cells.Item(7, 11).Value = ""
$ws3.Cells.Item(7, 12).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(7, 12).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Cells.Item(7, 13).Value = ""
$ws3.Cells.Item(7, 14).Value = ""
$ws3.Cells.Item(7, 15).Value = "'True"
$ws3.Cells.Item(7, 16).Value = ""
$ws3.Cells.Item(7, 17).Value = "'False"
$ws3.Cells.Item(7, 18).Value = ""

$ws3.Cells.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Cells.Item(2, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4a28272383e24ad73c8d92fe07e2d8956f1af8d3/e2e/17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md", "", "", "17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2, 10), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/629838ff0f85d56630f06315559779a76c260f11/e2e/17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md", "", "", "17c0919d-c3f6-4c5e-a05d-6cd3e4f8f9d0.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692db1b7948673fd1feb1023296c1e0e1ec9444a/e2e/23ea27e7-c26e-4c26-aa4c-af80ea04e4dc.md", "", "", "23ea27e7-c26e-4c26-aa4c-af80ea04e4dc.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1af3634fd0a2a83d98831a73295a81a69f4d53d3/e2e/38b81aa5-ef24-46c7-a1ab-2ffc89f9234c.md", "", "", "38b81aa5-ef24-46c7-a1ab-2ffc89f9234c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/108c5fd2cde90b86983be5bb4dee58dfd974048e/e2e/85e3d7a9-c5e2-4c10-8003-d5d93f2db402.md", "", "", "85e3d7a9-c5e2-4c10-8003-d5d93f2db402.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(6, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/692db1b7948673fd1feb1023296c1e0e1ec9444a/e2e/a7341649-f68e-4d0f-911b-33b6545d0ef9.md", "", "", "a7341649-f68e-4d0f-911b-33b6545d0ef9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(7, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/408e5bacc2e22a8dd0b7d42bca0f848975663a38/e2e/b013591f-7584-4a73-a059-a39a2aac7c99.md", "", "", "b013591f-7584-4a73-a059-a39a2aac7c99.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(8, 1), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/decaa4d912507b61ae7fd3fd8f9bb69ddc6d7150/e2e/b04048be-ba5b-487b-a279-7de58c31591f.md", "", "", "b04048be-ba5b-487b-a279-7de58c31591f.md") | Out-Null

$lo3 = $ws3.ListObjects.Item("de-de")
$lo3.Resize($ws3.Range("A1:R8"))
